$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Video file name
$ws.Range("B2").Value = "cut_Video_14.avi"

# Numeric parameters
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 1184
$ws.Range("G2").Value = 19
$ws.Range("H2").Value = 0.6
$ws.Range("M2").Value = 12
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 50
